$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 0.889714777469635
$ws.Range("B1").Value = 1.928547263145447
$ws.Range("C1").Value = 3.228466272354126
$ws.Range("D1").Value = 1.861977934837341
$ws.Range("E1").Value = 0.7267252802848816
